# Cassandra.docx update: expand the final paragraph with additional
# architecture/write-path details, splitting the original sentence into
# several runs and appending many new paragraphs (some with bold terms
# and spell-check proofErr markers), preserving the trailing _GoBack
# bookmark at the very end of the new content.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive original text instead of
# a hard-coded index, in case the surrounding content shifts.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ppText = $d.Paragraphs.Item($i).Range.Text
    if ($ppText -like "As a Cassandra is a distributed database system*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the target paragraph (Cassandra distributed database system sentence)."
}

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>As</w:t></w:r><w:r><w:t xml:space="preserve"> Cassandra is a distributed database system, data is distributed across</w:t></w:r><w:r><w:t xml:space="preserve"> all nodes and</w:t></w:r><w:r><w:t xml:space="preserve"> multiple servers which gives doorway to horizontal scalability.</w:t></w:r></w:p><w:p><w:r><w:t>Each node exchanges information with other nodes across the cluster.</w:t></w:r></w:p><w:p><w:r><w:t>After every write, commit log ensures data durability by capturing the write activity.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Data to be written to the node are first indexed and written to in-memory table called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>MemTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (can be called as “Write Back Cache”).</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">During, the process of write the data into the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MemTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, when the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MemTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is full, data is then written to the disk called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>SSTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>All writes are automatically partitioned and distributed across all nodes throughout the cluster</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Cassandra periodically consolidates SS Table by flushing out unwanted data using a process called </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>“Compaction”</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Tombstone</w:t></w:r><w:r><w:t xml:space="preserve"> is marker in a row that indicates a column was deleted.</w:t></w:r><w:r><w:t xml:space="preserve"> Hence, Cassandra is a row oriented database.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">An </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Authorized</w:t></w:r><w:r><w:t xml:space="preserve"> person can connect to any node and access data using </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>CQL</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In Cassandra, typically a cluster has </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">One </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Keyspace</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> per</w:t></w:r><w:r><w:t xml:space="preserve"> application.</w:t></w:r><w:r><w:t xml:space="preserve"> When a client gets connected to a Node, that Node serves as a coordinator for that client operation. A coordinator typically decides what nodes to be requested in the ring for fetching the data based on the partitioning of data and placement of replicas.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

[void]$r.InsertXML($xml)
